$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Remove legacy annotations (cell comments + hyperlinks) ---
$ws.Hyperlinks.Delete()
while ($ws.Comments.Count() -gt 0) {
  $ws.Comments.Item(1).Delete()
}

# --- Insert the new "GLAvailabilityChanges" column before the old StartDate column (E) ---
$ws.Columns("E:E").Insert()

# --- Drop the trailing blank formatting-only row (old row 4) ---
$ws.Rows("4:4").Delete()

# --- Reset formatting on the whole working area back to Normal so we can rebuild it cleanly ---
$ws.Range("A1:N3").Style = "Normal"

# --- Header row (row 1) ---
$ws.Range("A1").Value = "EmployeeID"
$ws.Range("B1").Value = "EmployeePassword"
$ws.Range("C1").Value = "ManagerID"
$ws.Range("D1").Value = "ManagerPassword"
$ws.Range("E1").Value = "GLAvailabilityChanges"
$ws.Range("F1").Value = "StartDate"
$ws.Range("G1").Value = "EndDate"
$ws.Range("H1").Value = "RepeatEvery"
$ws.Range("I1").Value = "DaysORWeeks"
$ws.Range("J1").Value = "Status"
$ws.Range("K1").Value = "StartTime"
$ws.Range("L1").Value = "EndTime"
$ws.Range("M1").Value = "Approve"
$ws.Range("N1").Value = "TestResult"

# --- Row 2 ---
$ws.Range("A2").Value = 10648995
$ws.Range("B2").Value = "WFMManagerPassword@05"
$ws.Range("C2").Value = "10648831"
$ws.Range("D2").Value = "WFMManagerPassword@06"
$ws.Range("E2").Value = "GL-Temporary Availability"
$ws.Range("F2").Value = "14/07/2025"
$ws.Range("G2").Value = "16/07/2025"
$ws.Range("H2").Value = 1
$ws.Range("I2").Value = "Day(s)"
$ws.Range("J2").Value = "Unavailable"
$ws.Range("K2").Value = 9
$ws.Range("L2").Value = 17
$ws.Range("M2").Value = "No"
$ws.Range("N2").Value = "Passed"

# --- Row 3 ---
$ws.Range("A3").Value = 10648995
$ws.Range("B3").Value = "WFMManagerPassword@05"
$ws.Range("C3").Value = "10648831"
$ws.Range("D3").Value = "WFMManagerPassword@06"
$ws.Range("E3").Value = "GL-Temporary Availability"
$ws.Range("F3").Value = "17/07/2025"
$ws.Range("G3").Value = "17/07/2025"
$ws.Range("H3").Value = 1
$ws.Range("I3").Value = "Day(s)"
$ws.Range("J3").Value = "Preferred Unavailable"
$ws.Range("K3").Value = 9
$ws.Range("L3").Value = 12
$ws.Range("M3").Value = "No"
$ws.Range("N3").Value = "Passed"

# --- Formatting: header row shaded with the theme accent tint ---
$ws.Range("A1:N1").Interior.ThemeColor = 10
$ws.Range("A1:N1").Interior.TintAndShade = 0.4

# --- Formatting: TestResult "Passed" cells shaded green ---
$ws.Range("N2").Interior.Color = 5287936
$ws.Range("N3").Interior.Color = 5287936

# --- Cosmetics: font size, row heights, view state ---
$ws.Range("A1:N3").Font.Size = 12
$ws.Rows("1:3").RowHeight = 15.5
$ws.Range("G7").Select()
